$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.811979666666667
$ws.Range("H2").Value = 8.435939000000001
$ws.Range("I2").Value = 0.01221198172659148
$ws.Range("J2").Value = 0.01221198172659148
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 1.153641840720111
$ws.Range("R2").Value = 10.382776566481
$ws.Range("S2").Value = 0.00004273729394302889
$ws.Range("T2").Value = 0.00004273729394302889

$ws.Range("G3").Value = 2.811979666666667
$ws.Range("H3").Value = 8.435939000000001
$ws.Range("I3").Value = 0.01221198172659148
$ws.Range("J3").Value = 0.01221198172659148
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 286.5355905464824
$ws.Range("R3").Value = 2578.820314918341
$ws.Range("S3").Value = 0.01061486791314755
$ws.Range("T3").Value = 0.01061486791314755

$ws.Range("G4").Value = 2.811979666666667
$ws.Range("H4").Value = 8.435939000000001
$ws.Range("I4").Value = 0.01221198172659148
$ws.Range("J4").Value = 0.01221198172659148
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 41.95852436327768
$ws.Range("R4").Value = 377.6267192694991
$ws.Range("S4").Value = 0.001554376519500899
$ws.Range("T4").Value = 0.001554376519500899

$ws.Range("I5").Value = 0.8095640809678946
$ws.Range("J5").Value = 0.8095640809678947
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 76.47792286775423
$ws.Range("R5").Value = 688.301305809788
$ws.Range("S5").Value = 0.002833166546483186
$ws.Range("T5").Value = 0.002833166546483186

$ws.Range("I6").Value = 0.8095640809678946
$ws.Range("J6").Value = 0.8095640809678947
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("S6").Value = 0.7036872457801674
$ws.Range("T6").Value = 0.7036872457801674

$ws.Range("I7").Value = 0.8095640809678946
$ws.Range("J7").Value = 0.8095640809678947
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("Q7").Value = 2781.539882340361
$ws.Range("R7").Value = 25033.85894106325
$ws.Range("S7").Value = 0.1030436686412441
$ws.Range("T7").Value = 0.1030436686412441

$ws.Range("I8").Value = 0.1782239373055139
$ws.Range("J8").Value = 0.1782239373055139
$ws.Range("M8").Value = 0.4102596666666667
$ws.Range("N8").Value = 1.230779
$ws.Range("O8").Value = 0.003499619873322347
$ws.Range("P8").Value = 0.003499619873322347
$ws.Range("Q8").Value = 16.83646403153489
$ws.Range("R8").Value = 151.528176283814
$ws.Range("S8").Value = 0.0006237160328961326
$ws.Range("T8").Value = 0.0006237160328961325

$ws.Range("I9").Value = 0.1782239373055139
$ws.Range("J9").Value = 0.1782239373055139
$ws.Range("O9").Value = 0.8692174743460166
$ws.Range("P9").Value = 0.8692174743460165
$ws.Range("S9").Value = 0.1549153606527016
$ws.Range("T9").Value = 0.1549153606527016

$ws.Range("I10").Value = 0.1782239373055139
$ws.Range("J10").Value = 0.1782239373055139
$ws.Range("N10").Value = 44.764041
$ws.Range("O10").Value = 0.1272829057806611
$ws.Range("P10").Value = 0.1272829057806611
$ws.Range("Q10").Value = 612.3505245073674
$ws.Range("R10").Value = 5511.154720566306
$ws.Range("S10").Value = 0.02268486061991619
$ws.Range("T10").Value = 0.02268486061991619
